$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.626.94'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '1.562.68'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.521'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.67%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.90'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0588'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0901'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').Value = '1.785.34'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.558.63'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '28.667.40'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.515'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').Value = '0.0₃0684'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  +3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0459'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.99%  '
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').Value = '1.404.98'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.24%  '
$ws.Range('E38').Value = '  -2.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0162'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.514'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.771'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0457'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.51%  '
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('D47').Value = '1.697.37'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.852'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.85%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0106'
$ws.Range('E49').Value = '  +4.23%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '84.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.09'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.92%  '
